# Payslip template employees workbook update.
#
# The underlying change (per the commit's folder reorg / Google Drive
# documentation-folder addition) only touches the payslip numbers for the
# "Duncan" employee row: Hours Worked goes from 8 to 13, and the two
# Total Cost / Currency Pay cells derived from it go from 640 to 1040
# (13 hrs * 80 rate).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Hours Worked (row 9)
$ws.Range("G9").Value = 13

# Total Cost (row 9) and Currency Pay (row 11) recompute to the new total
$ws.Range("I9").Value = 1040
$ws.Range("I11").Value = 1040
